# Update res_bus/vm_pu.xlsx results for "case with 380 kV done"
# Bus 1 (slack) voltage setpoint lowered from 1.05 pu to 1.02 pu,
# and all dependent bus voltage results (rows 2-25) recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033169590446045
$ws.Range("D2").Value = 1.03753826882599
$ws.Range("E2").Value = 1.052123864072763
$ws.Range("F2").Value = 1.057553981660915
$ws.Range("I2").Value = 1.037922478219318
$ws.Range("J2").Value = 1.038295790881637
$ws.Range("K2").Value = 1.040328867398276
$ws.Range("L2").Value = 1.054873466286093
$ws.Range("M2").Value = 1.060288634187176
$ws.Range("N2").Value = 1.039770289998881

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034039402005496
$ws.Range("D3").Value = 1.038201672324784
$ws.Range("E3").Value = 1.053360078097911
$ws.Range("F3").Value = 1.058855526115014
$ws.Range("I3").Value = 1.038146384996785
$ws.Range("J3").Value = 1.03880860791844
$ws.Range("K3").Value = 1.040802486836172
$ws.Range("L3").Value = 1.055921376181862
$ws.Range("M3").Value = 1.061402799449501
$ws.Range("N3").Value = 1.040283835294698

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03460212914186
$ws.Range("D4").Value = 1.038630716912878
$ws.Range("E4").Value = 1.054160790830322
$ws.Range("F4").Value = 1.059698482950441
$ws.Range("I4").Value = 1.038289788739589
$ws.Range("J4").Value = 1.039139692872616
$ws.Range("K4").Value = 1.041108046951252
$ws.Range("L4").Value = 1.056599679961033
$ws.Range("M4").Value = 1.062123957806312
$ws.Range("N4").Value = 1.040615390427484

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03483867516092
$ws.Range("D5").Value = 1.038811033037115
$ws.Range("E5").Value = 1.054497603011185
$ws.Range("F5").Value = 1.060053047616094
$ws.Range("I5").Value = 1.03834972115403
$ws.Range("J5").Value = 1.039278702746202
$ws.Range("K5").Value = 1.041236287436277
$ws.Range("L5").Value = 1.056884896320435
$ws.Range("M5").Value = 1.062427186092869
$ws.Range("N5").Value = 1.040754597711035

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034878390821875
$ws.Range("D6").Value = 1.038841305695838
$ws.Range("E6").Value = 1.054554166582321
$ws.Range("F6").Value = 1.060112591501363
$ws.Range("I6").Value = 1.038359763267452
$ws.Range("J6").Value = 1.039302032645597
$ws.Range("K6").Value = 1.041257806832531
$ws.Range("L6").Value = 1.056932788798887
$ws.Range("M6").Value = 1.06247810260909
$ws.Range("N6").Value = 1.040777960741563

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034605289978423
$ws.Range("D7").Value = 1.038633126520212
$ws.Range("E7").Value = 1.054165290573164
$ws.Range("F7").Value = 1.059703219930085
$ws.Range("I7").Value = 1.038290590953045
$ws.Range("J7").Value = 1.039141551030492
$ws.Range("K7").Value = 1.041109761360421
$ws.Range("L7").Value = 1.056603490809045
$ws.Range("M7").Value = 1.062128009346896
$ws.Range("N7").Value = 1.040617251224157

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033463567409383
$ws.Range("D8").Value = 1.0377625144099
$ws.Range("E8").Value = 1.05254148350357
$ws.Range("F8").Value = 1.057993685843297
$ws.Range("I8").Value = 1.037998454503741
$ws.Range("J8").Value = 1.038469253059385
$ws.Range("K8").Value = 1.040489115908944
$ws.Range("L8").Value = 1.055227564045886
$ws.Range("M8").Value = 1.060665126764093
$ws.Range("N8").Value = 1.039943998512817

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031450979752627
$ws.Range("D9").Value = 1.036226734386599
$ws.Range("E9").Value = 1.049686201600812
$ws.Range("F9").Value = 1.054987118572494
$ws.Range("I9").Value = 1.037472366314304
$ws.Range("J9").Value = 1.037278922415161
$ws.Range("K9").Value = 1.039388572704939
$ws.Range("L9").Value = 1.052804770211395
$ws.Range("M9").Value = 1.058088977039187
$ws.Range("N9").Value = 1.038751977462546

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030108807741054
$ws.Range("D10").Value = 1.03520183592457
$ws.Range("E10").Value = 1.047786691431908
$ws.Range("F10").Value = 1.05298659981944
$ws.Range("I10").Value = 1.037114065485303
$ws.Range("J10").Value = 1.036481600778681
$ws.Range("K10").Value = 1.038650288671098
$ws.Range("L10").Value = 1.051190699454978
$ws.Range("M10").Value = 1.056372585256007
$ws.Range("N10").Value = 1.037953523537888

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029527534550621
$ws.Range("D11").Value = 1.034757808860507
$ws.Range("E11").Value = 1.046965117141634
$ws.Range("F11").Value = 1.052121251331716
$ws.Range("I11").Value = 1.036957125766216
$ws.Range("J11").Value = 1.036135464874189
$ws.Range("K11").Value = 1.038329523542049
$ws.Range("L11").Value = 1.050492041555515
$ws.Range("M11").Value = 1.05562960186808
$ws.Range("N11").Value = 1.03760689608071

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029311608775864
$ws.Range("D12").Value = 1.034592842504421
$ws.Range("E12").Value = 1.046660085482695
$ws.Range("F12").Value = 1.051799954414365
$ws.Range("I12").Value = 1.036898562347717
$ws.Range("J12").Value = 1.036006761309214
$ws.Range("K12").Value = 1.038210214957252
$ws.Range("L12").Value = 1.050232564456249
$ws.Range("M12").Value = 1.055353657367177
$ws.Range("N12").Value = 1.037478009741909

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029357926258532
$ws.Range("D13").Value = 1.034628229912436
$ws.Range("E13").Value = 1.046725509606243
$ws.Range("F13").Value = 1.051868867743888
$ws.Range("I13").Value = 1.036911136567485
$ws.Range("J13").Value = 1.036034374684078
$ws.Range("K13").Value = 1.038235814380746
$ws.Range("L13").Value = 1.05028822158711
$ws.Range("M13").Value = 1.055412846944154
$ws.Range("N13").Value = 1.037505662330933

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029509686362038
$ws.Range("D14").Value = 1.034744173396861
$ws.Range("E14").Value = 1.04693990033112
$ws.Range("F14").Value = 1.052094690143785
$ws.Range("I14").Value = 1.036952290381527
$ws.Range("J14").Value = 1.03612482891337
$ws.Range("K14").Value = 1.038319664761144
$ws.Range("L14").Value = 1.050470592379876
$ws.Range("M14").Value = 1.055606791551695
$ws.Range("N14").Value = 1.037596245015606

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029603188795903
$ws.Range("D15").Value = 1.034815605405314
$ws.Range("E15").Value = 1.047072011721872
$ws.Range("F15").Value = 1.052233844227162
$ws.Range("I15").Value = 1.036977610968111
$ws.Range("J15").Value = 1.036180543103352
$ws.Range("K15").Value = 1.038371306281966
$ws.Range("L15").Value = 1.050582961716072
$ws.Range("M15").Value = 1.055726291508387
$ws.Range("N15").Value = 1.037652038326129

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03014738278914
$ws.Range("D16").Value = 1.035231299583963
$ws.Range("E16").Value = 1.047841235993935
$ws.Range("F16").Value = 1.053044048751612
$ws.Range("I16").Value = 1.037124443311014
$ws.Range("J16").Value = 1.036504553961441
$ws.Range("K16").Value = 1.038671554025636
$ws.Range("L16").Value = 1.051237072136562
$ws.Range("M16").Value = 1.056421899259875
$ws.Range("N16").Value = 1.0379765093168

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03048871390395
$ws.Range("D17").Value = 1.035491990218215
$ws.Range("E17").Value = 1.048323996932262
$ws.Range("F17").Value = 1.053552505634165
$ws.Range("I17").Value = 1.037216067561797
$ws.Range("J17").Value = 1.036707559381063
$ws.Range("K17").Value = 1.038859601916107
$ws.Range("L17").Value = 1.051647443090821
$ws.Range("M17").Value = 1.056858295114575
$ws.Range("N17").Value = 1.038179803027403

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030687796509866
$ws.Range("D18").Value = 1.035644023501274
$ws.Range("E18").Value = 1.048605672682185
$ws.Range("F18").Value = 1.053849165844391
$ws.Range("I18").Value = 1.037269337404428
$ws.Range("J18").Value = 1.03682588306681
$ws.Range("K18").Value = 1.038969182508142
$ws.Range("L18").Value = 1.051886829289999
$ws.Range("M18").Value = 1.057112859284333
$ws.Range("N18").Value = 1.038298294746356

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030755676810198
$ws.Range("D19").Value = 1.035695858963787
$ws.Range("E19").Value = 1.048701732054707
$ws.Range("F19").Value = 1.053950333929895
$ws.Range("I19").Value = 1.037287471683427
$ws.Range("J19").Value = 1.036866213793858
$ws.Range("K19").Value = 1.039006528904084
$ws.Range("L19").Value = 1.05196845786854
$ws.Range("M19").Value = 1.057199662832987
$ws.Range("N19").Value = 1.038338682747662

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030452093341838
$ws.Range("D20").Value = 1.035464022967801
$ws.Range("E20").Value = 1.048272192009725
$ws.Range("F20").Value = 1.053497944147397
$ws.Range("I20").Value = 1.037206255041661
$ws.Range("J20").Value = 1.036685787706434
$ws.Range("K20").Value = 1.038839436967857
$ws.Range("L20").Value = 1.051603411722379
$ws.Range("M20").Value = 1.056811471718283
$ws.Range("N20").Value = 1.038158000434499

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029464997204411
$ws.Range("D21").Value = 1.03471003188676
$ws.Range("E21").Value = 1.046876763819668
$ws.Range("F21").Value = 1.052028187443126
$ws.Range("I21").Value = 1.036940179032275
$ws.Range("J21").Value = 1.036098196072565
$ws.Range("K21").Value = 1.038294977384794
$ws.Range("L21").Value = 1.050416887755999
$ws.Range("M21").Value = 1.055549678808765
$ws.Range("N21").Value = 1.037569574353112

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02884428400443
$ws.Range("D22").Value = 1.034235766147161
$ws.Range("E22").Value = 1.046000198164197
$ws.Range("F22").Value = 1.051104856239656
$ws.Range("I22").Value = 1.036771330207078
$ws.Range("J22").Value = 1.035727983024363
$ws.Range("K22").Value = 1.037951716281408
$ws.Range("L22").Value = 1.04967107809214
$ws.Range("M22").Value = 1.054756527570401
$ws.Range("N22").Value = 1.037198835559918

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029173343841818
$ws.Range("D23").Value = 1.034487202163774
$ws.Range("E23").Value = 1.046464807315678
$ws.Range("F23").Value = 1.051594259764852
$ws.Range("I23").Value = 1.036860987556874
$ws.Range("J23").Value = 1.035924312872108
$ws.Range("K23").Value = 1.038133774170713
$ws.Range("L23").Value = 1.0500664269082
$ws.Range("M23").Value = 1.055176974651079
$ws.Range("N23").Value = 1.037395444218566

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030468640628596
$ws.Range("D24").Value = 1.035476660236033
$ws.Range("E24").Value = 1.048295600147447
$ws.Range("F24").Value = 1.053522597868785
$ws.Range("I24").Value = 1.037210689431775
$ws.Range("J24").Value = 1.03669562565458
$ws.Range("K24").Value = 1.038848548963069
$ws.Range("L24").Value = 1.051623307528659
$ws.Range("M24").Value = 1.056832629127672
$ws.Range("N24").Value = 1.03816785235366

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031971363712651
$ws.Range("D25").Value = 1.036623959237798
$ws.Range("E25").Value = 1.050423648135573
$ws.Range("F25").Value = 1.055763702626669
$ws.Range("I25").Value = 1.037609709862073
$ws.Range("J25").Value = 1.037587318062345
$ws.Range("K25").Value = 1.039673901301249
$ws.Range("L25").Value = 1.053430917687435
$ws.Range("M25").Value = 1.058754785650758
$ws.Range("N25").Value = 1.039060811066923
